$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "27.070.76"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.825.82"
$ws.Range("E3").Value = "  -1.65%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.96%  "
Set-TextValue $ws.Range("D5") "311.27"
$ws.Range("E5").Value = "  -2.59%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.91%  "
Set-TextValue $ws.Range("D7") "0.4230"
$ws.Range("E7").Value = "  -1.69%  "
Set-TextValue $ws.Range("D8") "0.3678"
$ws.Range("E8").Value = "  -1.95%  "
Set-TextValue $ws.Range("D9") "0.07231"
$ws.Range("E9").Value = "  -1.50%  "
Set-TextValue $ws.Range("D10") "0.8433"
$ws.Range("E10").Value = "  -3.97%  "
Set-TextValue $ws.Range("D11") "20.77"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").Value = "1.819.78"
$ws.Range("E12").Value = "  -2.05%  "
Set-TextValue $ws.Range("D13") "6.687"
$ws.Range("E13").Value = "  -0.92%  "
Set-TextValue $ws.Range("D14") "0.07074"
$ws.Range("E14").Value = "  -0.87%  "
Set-TextValue $ws.Range("D15") "5.299"
$ws.Range("E15").Value = "  -2.68%  "
Set-TextValue $ws.Range("D16") "90.15"
$ws.Range("E16").Value = "  +1.30%  "
Set-TextValue $ws.Range("D17") "1.002"
$ws.Range("E17").Value = "  -1.08%  "
Set-TextValue $ws.Range("D18") "0.000008739"
$ws.Range("E18").Value = "  -2.90%  "
Set-TextValue $ws.Range("D19") "1.000"
$ws.Range("E19").Value = "  -0.91%  "
Set-TextValue $ws.Range("D20") "14.92"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("D21").Value = "27.166.30"
$ws.Range("E21").Value = "  -2.18%  "
Set-TextValue $ws.Range("D22") "5.149"
$ws.Range("E22").Value = "  -1.20%  "
Set-TextValue $ws.Range("D23") "10.86"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("D24").Value = "2.050.75"
$ws.Range("E24").Value = "  -1.23%  "
Set-TextValue $ws.Range("D25") "1.992"
$ws.Range("E25").Value = "  +0.44%  "
Set-TextValue $ws.Range("D26") "152.04"
$ws.Range("E26").Value = "  -2.05%  "
Set-TextValue $ws.Range("D27") "2.291"
$ws.Range("E27").Value = "  +4.95%  "
Set-TextValue $ws.Range("D28") "18.27"
$ws.Range("E28").Value = "  -2.03%  "
Set-TextValue $ws.Range("D29") "5.269"
$ws.Range("E29").Value = "  -1.95%  "
Set-TextValue $ws.Range("D30") "117.24"
$ws.Range("E30").Value = "  -1.54%  "
Set-TextValue $ws.Range("D31") "0.08712"
$ws.Range("E31").Value = "  -2.61%  "
Set-TextValue $ws.Range("D32") "1.181"
$ws.Range("E32").Value = "  -3.88%  "
Set-TextValue $ws.Range("D33") "0.7377"
$ws.Range("E33").Value = "  -5.44%  "
Set-TextValue $ws.Range("D34") "4.435"
$ws.Range("E34").Value = "  -2.67%  "
Set-TextValue $ws.Range("D35") "2.886"
$ws.Range("E35").Value = "  -1.50%  "
Set-TextValue $ws.Range("D36") "0.9999"
$ws.Range("E36").Value = "  -1.09%  "
Set-TextValue $ws.Range("D37") "1.089"
$ws.Range("E37").Value = "  -3.64%  "
Set-TextValue $ws.Range("D38") "0.01953"
$ws.Range("E38").Value = "  -1.56%  "
Set-TextValue $ws.Range("D39") "0.05266"
$ws.Range("E39").Value = "  -1.66%  "
Set-TextValue $ws.Range("D40") "7.342"
$ws.Range("E40").Value = "  +1.11%  "
Set-TextValue $ws.Range("D41") "2.878"
$ws.Range("E41").Value = "  -0.85%  "
Set-TextValue $ws.Range("D42") "0.1689"
$ws.Range("E42").Value = "  -0.45%  "
Set-TextValue $ws.Range("D43") "0.5064"
$ws.Range("E43").Value = "  -1.35%  "
Set-TextValue $ws.Range("D44") "8.583"
$ws.Range("E44").Value = "  -2.95%  "
Set-TextValue $ws.Range("D45") "10.51"
$ws.Range("E45").Value = "  -2.06%  "
Set-TextValue $ws.Range("D46") "106.24"
$ws.Range("E46").Value = "  -1.57%  "
Set-TextValue $ws.Range("D47") "0.4725"
$ws.Range("E47").Value = "  -1.16%  "
Set-TextValue $ws.Range("D48") "1.924"
$ws.Range("E48").Value = "  +4.41%  "
Set-TextValue $ws.Range("D49") "0.9996"
$ws.Range("E49").Value = "  -1.12%  "
Set-TextValue $ws.Range("D50") "0.06331"
$ws.Range("E50").Value = "  -2.26%  "
Set-TextValue $ws.Range("D51") "1.655"
$ws.Range("E51").Value = "  -1.93%  "
